# Updated test data for German, Czech market
# Added test data for Belgium market
#
# Removes the obsolete "PR1DS / PR8AS / ZXF / ZXFEV" product rows from the
# Germany and Belgium sheets (rows 18-21), and refreshes the Jira
# reference strings (B4) on the Germany and Belgium sheets with the extra
# ticket numbers that were added ("/T1746" and "/T2267" respectively).
# Czech's rows/reference were already up to date.

$wb = $excel.ActiveWorkbook

# --- Belgium sheet -----------------------------------------------------
$wsBelgium = $wb.Worksheets.Item("Belgium")
$wsBelgium.Rows("18:21").Delete()
$wsBelgium.Range("B4").Value = "NGC-3478/T2265/T2264/T2268/T2267"
$wsBelgium.Range("A15").Select()

# --- Germany sheet ------------------------------------------------------
$wsGermany = $wb.Worksheets.Item("Germany")
$wsGermany.Rows("18:21").Delete()
$wsGermany.Range("B4").Value = "NGC-3475/T1730/T1746"

# --- Czech sheet ---------------------------------------------------------
$wsCzech = $wb.Worksheets.Item("Czech")
$wsCzech.Range("A13").Select()

# Leave Germany as the active sheet/selection, matching the saved view state.
$wsGermany.Select()
$wsGermany.Range("A11").Select()
